# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Hyperion_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1028.8823
$ws.Range("J2").Value = 944.75
$ws.Range("L2").Value = 944.75
$ws.Range("N2").Value = -1170.75
$ws.Range("H40").Value = 5081.4
$ws.Range("I40").Value = 3790.4
$ws.Range("J40").Value = 6372.4
$ws.Range("K40").Value = 3790.4
$ws.Range("L40").Value = 6372.4
$ws.Range("M40").Value = -3615.4
$ws.Range("N40").Value = -6722.4
$ws.Range("H76").Value = 7942336.5
$ws.Range("I76").Value = 13893589
$ws.Range("K76").Value = 13893589
$ws.Range("M76").Value = -13893274
$ws.Range("H79").Value = 7942336.5
$ws.Range("I79").Value = 13893589
$ws.Range("K79").Value = 13893589
$ws.Range("M79").Value = -13892497
$ws.Range("H82").Value = 1772.1111
$ws.Range("I82").Value = 1772.1111
$ws.Range("K82").Value = 5316.3333
$ws.Range("M82").Value = -4910.3333
$ws.Range("H85").Value = 1772.1111
$ws.Range("I85").Value = 1772.1111
$ws.Range("K85").Value = 5316.3333
$ws.Range("M85").Value = -3912.3333
$ws.Range("H107").Value = 847.1111
$ws.Range("I107").Value = 822
$ws.Range("K107").Value = 822
$ws.Range("M107").Value = 1098
$ws.Range("H132").Value = 4000.0344
$ws.Range("I132").Value = 4152.2607
$ws.Range("K132").Value = 12456.7821
$ws.Range("M132").Value = -9926.7821

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1750
$ws.Range("I35").Value = 1750
$ws.Range("K35").Value = 1750
$ws.Range("M35").Value = -1344
$ws.Range("H43").Value = 34949.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 34949.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 34949.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -35575.5
$ws.Range("H45").Value = 42840.16
$ws.Range("I45").Value = 63723.25
$ws.Range("K45").Value = 63723.25
$ws.Range("M45").Value = -63346.25
$ws.Range("H61").Value = 2936
$ws.Range("I61").Value = 2641.8572
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 2641.8572
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -2429.8572
$ws.Range("N61").Value = -5419
$ws.Range("H110").Value = 2388.4375
$ws.Range("I110").Value = 2392.4614
$ws.Range("K110").Value = 2392.4614
$ws.Range("M110").Value = -347.4614000000001
$ws.Range("H132").Value = 2927.8333
$ws.Range("I132").Value = 2547.5
$ws.Range("J132").Value = 3003.9
$ws.Range("K132").Value = 7642.5
$ws.Range("L132").Value = 9011.700000000001
$ws.Range("M132").Value = -5112.5
$ws.Range("N132").Value = -14071.7
$ws.Range("H136").Value = 2936
$ws.Range("I136").Value = 2641.8572
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 7925.571599999999
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -5375.571599999999
$ws.Range("N136").Value = -20085

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 10048
$ws.Range("I39").Value = 10048
$ws.Range("K39").Value = 10048
$ws.Range("M39").Value = -9659
$ws.Range("H86").Value = 5782.0415
$ws.Range("I86").Value = 6154.273
$ws.Range("J86").Value = 1687.5
$ws.Range("K86").Value = 6154.273
$ws.Range("L86").Value = 1687.5
$ws.Range("M86").Value = -5031.273
$ws.Range("N86").Value = -3933.5
$ws.Range("H89").Value = 5782.0415
$ws.Range("I89").Value = 6154.273
$ws.Range("J89").Value = 1687.5
$ws.Range("K89").Value = 30771.365
$ws.Range("L89").Value = 8437.5
$ws.Range("M89").Value = -25155.365
$ws.Range("N89").Value = -19669.5
$ws.Range("H94").Value = 4049.8276
$ws.Range("I94").Value = 1033.1578
$ws.Range("J94").Value = 9781.5
$ws.Range("K94").Value = 1033.1578
$ws.Range("L94").Value = 9781.5
$ws.Range("M94").Value = -582.1578
$ws.Range("N94").Value = -10683.5
$ws.Range("H105").Value = 1591.8125
$ws.Range("I105").Value = 1677.0714
$ws.Range("J105").Value = 995
$ws.Range("K105").Value = 1677.0714
$ws.Range("L105").Value = 995
$ws.Range("M105").Value = 69.92859999999996
$ws.Range("N105").Value = -4489
$ws.Range("H134").Value = 3283.6316
$ws.Range("I134").Value = 1582.0667
$ws.Range("J134").Value = 9664.5
$ws.Range("K134").Value = 4746.2001
$ws.Range("L134").Value = 28993.5
$ws.Range("M134").Value = -2211.2001
$ws.Range("N134").Value = -34063.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 27250
$ws.Range("I23").Value = 8000
$ws.Range("J23").Value = 30000
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = -7760
$ws.Range("N23").Value = -30480
$ws.Range("H27").Value = 27250
$ws.Range("I27").Value = 8000
$ws.Range("J27").Value = 30000
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = -7808
$ws.Range("N27").Value = -30384
$ws.Range("H134").Value = 3744.8462
$ws.Range("I134").Value = 3405.9333
$ws.Range("J134").Value = 4207
$ws.Range("K134").Value = 10217.7999
$ws.Range("L134").Value = 12621
$ws.Range("M134").Value = -7682.7999
$ws.Range("N134").Value = -17691

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 61622.824
$ws.Range("J55").Value = 94772.55
$ws.Range("L55").Value = 284317.65
$ws.Range("N55").Value = -284671.65
$ws.Range("H59").Value = 2500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 7500
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -8580
$ws.Range("H114").Value = 166666930
$ws.Range("I114").Value = 166666930
$ws.Range("K114").Value = 500000790
$ws.Range("M114").Value = -499997536
$ws.Range("H131").Value = 13023294
$ws.Range("J131").Value = 16670219
$ws.Range("L131").Value = 50010657
$ws.Range("N131").Value = -50020737

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 85085.164
$ws.Range("J51").Value = 85085.164
$ws.Range("L51").Value = 85085.164
$ws.Range("N51").Value = -86103.164
$ws.Range("H132").Value = 3969.6875
$ws.Range("I132").Value = 3247.75
$ws.Range("J132").Value = 4691.625
$ws.Range("K132").Value = 9743.25
$ws.Range("L132").Value = 14074.875
$ws.Range("M132").Value = -7213.25
$ws.Range("N132").Value = -19134.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5932.1665
$ws.Range("I40").Value = 4626
$ws.Range("J40").Value = 7984.7144
$ws.Range("K40").Value = 4626
$ws.Range("L40").Value = 7984.7144
$ws.Range("M40").Value = -4490
$ws.Range("N40").Value = -8256.714400000001
$ws.Range("H55").Value = 2678.2354
$ws.Range("I55").Value = 2424.8572
$ws.Range("K55").Value = 2424.8572
$ws.Range("M55").Value = -2251.8572
$ws.Range("H82").Value = 1945.2
$ws.Range("I82").Value = 1445.3334
$ws.Range("J82").Value = 2695
$ws.Range("K82").Value = 1445.3334
$ws.Range("L82").Value = 2695
$ws.Range("M82").Value = -1084.3334
$ws.Range("N82").Value = -3417
$ws.Range("H85").Value = 1945.2
$ws.Range("I85").Value = 1445.3334
$ws.Range("J85").Value = 2695
$ws.Range("K85").Value = 1445.3334
$ws.Range("L85").Value = 2695
$ws.Range("M85").Value = -197.3334
$ws.Range("N85").Value = -5191
$ws.Range("H107").Value = 3633.5
$ws.Range("I107").Value = 3633.5
$ws.Range("K107").Value = 3633.5
$ws.Range("M107").Value = -1713.5
$ws.Range("H136").Value = 33696.91
$ws.Range("I136").Value = 52390.25
$ws.Range("J136").Value = 4937.923
$ws.Range("K136").Value = 157170.75
$ws.Range("L136").Value = 14813.769
$ws.Range("M136").Value = -154620.75
$ws.Range("N136").Value = -19913.769

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 8000
$ws.Range("I58").Value = 8000
$ws.Range("K58").Value = 8000
$ws.Range("M58").Value = -7692
$ws.Range("H62").Value = 6765.8
$ws.Range("I62").Value = 3576.8845
$ws.Range("K62").Value = 3576.8845
$ws.Range("M62").Value = -2952.8845
$ws.Range("H65").Value = 6765.8
$ws.Range("I65").Value = 3576.8845
$ws.Range("K65").Value = 17884.4225
$ws.Range("M65").Value = -14764.4225
$ws.Range("H107").Value = 2507.7144
$ws.Range("I107").Value = 2932.8823
$ws.Range("K107").Value = 8798.6469
$ws.Range("M107").Value = -6878.6469
$ws.Range("H122").Value = 1662
$ws.Range("I122").Value = 1077.2778
$ws.Range("K122").Value = 3231.8334
$ws.Range("M122").Value = -781.8334000000004
$ws.Range("H132").Value = 14587.176
$ws.Range("I132").Value = 2595.5156
$ws.Range("K132").Value = 7786.5468
$ws.Range("M132").Value = -5256.5468
$ws.Range("H136").Value = 3699.3333
$ws.Range("I136").Value = 3189.2
$ws.Range("K136").Value = 9567.599999999999
$ws.Range("M136").Value = -7017.599999999999
